$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking values (e.g. "45.13") are not silently converted to
# numbers by Excel and lose formatting such as trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.814.45"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "2.213.88"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D5").Value = "261.33"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").Value = "86.93"
$ws.Range("E6").Value = "  +13.89%  "
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "45.13"
$ws.Range("E10").Value = "  +8.10%  "
$ws.Range("D11").Value = "0.0919"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "7.46"
$ws.Range("E12").Value = "  +7.57%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "2.547.48"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "14.48"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "2.209.39"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "43.754.58"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "69.84"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("D22").Value = "2.34"
$ws.Range("E22").Value = "  +6.28%  "
$ws.Range("D23").Value = "230.85"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "8.89"
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +5.57%  "
$ws.Range("D27").Value = "10.66"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").Value = "39.92"
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "174.83"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "20.47"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "0.0877"
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("D35").Value = "0.123"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("E36").Value = "  +4.75%  "
$ws.Range("D37").Value = "4.50"
$ws.Range("E37").Value = "  +4.67%  "
$ws.Range("D38").Value = "0.0358"
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").Value = "2.98"
$ws.Range("E39").Value = "  +8.60%  "
$ws.Range("D40").Value = "12.57"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "2.10"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "63.93"
$ws.Range("E42").Value = "  +5.90%  "
$ws.Range("D43").Value = "5.53"
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "100.72"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "0.0979"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("B48").Value = "WOONetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D48").Value = "0.453"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "1.18"
$ws.Range("E49").Value = "  +4.44%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "1.12"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("E51").Value = "  +5.15%  "
